$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 6.858228537776629
$ws.Range("E2").Value = 6.097510900342726
$ws.Range("H2").Value = 6.845950708523831
$ws.Range("I2").Value = 6.099200962514638
